# Swap the presentation's theme colour palette from the "Integral" scheme
# to the default "Office Theme" scheme (what the XML diff shows happening
# to ppt/theme/theme1.xml, the theme used by the slide master).
#
# Helper: convert an RRGGBB hex string into the decimal BBGGRR value that
# the PowerPoint COM RGB properties use.
function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$colorScheme = $design.SlideMaster.Theme.ThemeColorScheme

# Index order matches the OOXML <a:clrScheme> child order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeTheme = @{
    1  = "000000"
    2  = "FFFFFF"
    3  = "44546A"
    4  = "E7E6E6"
    5  = "5B9BD5"
    6  = "ED7D31"
    7  = "A5A5A5"
    8  = "FFC000"
    9  = "4472C4"
    10 = "70AD47"
    11 = "0563C1"
    12 = "954F72"
}

foreach ($idx in $officeTheme.Keys) {
    $colorScheme.Item($idx).RGB = HexToComRgb $officeTheme[$idx]
}
